$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$theme = $nm.Theme
$fs = $theme.ThemeFontScheme
$major = $fs.MajorFont
$major.Name = "TESTFONTNAME"
Write-Output "done"
